# Update CV results workbook with refreshed band / frequency-domain feature
# cross-validation numbers (commit: "band freq domain features upd").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "CV_Summary": per-frequency aggregate stats (rows 2-5)
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("CV_Summary")

$summaryData = @(
    @("30hz", 0.9375,  0.05929270612815711, 1,       0.84375, 15.81138563417568, 160, 246),
    @("10hz", 0.93125, 0.02338535866733714, 0.96875, 0.90625, 39.82190801638708, 160, 246),
    @("20hz", 0.8875,  0.0318688719599549,  0.9375,  0.84375, 27.84848245116334, 160, 246),
    @("40hz", 0.875,   0.05590169943749474, 0.9375,  0.78125, 15.65247304249903, 160, 246)
)

$row = 2
foreach ($rec in $summaryData) {
    $wsSummary.Cells.Item($row, 1).Value = $rec[0]
    $wsSummary.Cells.Item($row, 2).Value = $rec[1]
    $wsSummary.Cells.Item($row, 3).Value = $rec[2]
    $wsSummary.Cells.Item($row, 4).Value = $rec[3]
    $wsSummary.Cells.Item($row, 5).Value = $rec[4]
    $wsSummary.Cells.Item($row, 6).Value = $rec[5]
    $wsSummary.Cells.Item($row, 7).Value = $rec[6]
    $wsSummary.Cells.Item($row, 8).Value = $rec[7]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Sheet "CV_Scores_Detail": per-fold accuracy (column C), rows 2-21
# ---------------------------------------------------------------------------
$wsDetail = $wb.Worksheets.Item("CV_Scores_Detail")

$detailAccuracy = @(
    0.9375,
    0.90625,
    0.90625,
    0.9375,
    0.96875,
    0.875,
    0.9375,
    0.90625,
    0.875,
    0.84375,
    1,
    0.9375,
    1,
    0.90625,
    0.84375,
    0.90625,
    0.78125,
    0.84375,
    0.9375,
    0.90625
)

$row = 2
foreach ($acc in $detailAccuracy) {
    $wsDetail.Cells.Item($row, 3).Value = $acc
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# Sheet "Analysis_Info": summary labels + refreshed analysis timestamp
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Analysis_Info")

$wsInfo.Range("B3").Value = "30hz"
$wsInfo.Range("B4").Value = "40hz"
$wsInfo.Range("B5").Value = "10hz"
$wsInfo.Range("B6").Value = "2025-10-03 16:22:27"
